# Apply edit: Software Design 3.2 Functions section insertion
$d = $word.ActiveDocument

# --- Step 1: locate the "Functions" heading paragraph (Heading3, exact text match) precisely ---
$targetIndex = 0
$i = 0
$target = $null
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -eq "Functions`r") {
        $targetIndex = $i
        $target = $p
        break
    }
}
if ($target -eq $null) {
    Write-Host "ERROR: could not find 'Functions' heading paragraph"
} else {

    # --- Step 2: grab a ListTemplate reference from the document's existing numId=7 bullet list
    #     (used later so the new bullet items re-join that exact same list/numbering) ---
    $listSource = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match "reporting information of all listings") {
            $listSource = $p
            break
        }
    }
    $listTemplate = $listSource.Range.ListFormat.ListTemplate

    # --- Step 3: insert all seven new paragraphs in a single Find/Replace pass, scoped tightly
    #     to the "Functions" paragraph's own range so nothing else in the document can match.
    #     Using the Find/Replace "^p" paragraph-mark code (rather than Range.InsertParagraphAfter,
    #     which mis-targets in this runtime) keeps formatting clean and avoids corrupting
    #     unrelated parts of the document. ---
    $rng = $target.Range
    $ok = $rng.Find.Execute("Functions", $true, $false, $false, $false, $false, $true, 1, $false, "Functions^pThe following functions will be provided by the software:^pSearching the dataset for keywords, ordering the data accordingly and displaying the total result per search.^pUsing 2 fields from the dataset, displaying select data as a viewable graph.^pSearching the dataset using search parameters such as selected fields and periods.^pSearch function: When searching the dataset for keywords, the software uses the user’s query as a string variable, and searches the given dataset for matching strings, counting the total matches as a results variable. Once the results are collated, the display screen is ordered ascending by the related results and the results variable is displayed. Depending on the character size of the largest datatype, the search function will be assigned the same limit to account for all possible inputs without causing extra bloat by allowing too many characters. The setback for this function however is the possible hardware strain it may cause due to searching the entire dataset, hence care should be taken to minimise the strain, and an efficient coding structure should be prioritised utilising coding functions such as arrays to store the result ids to call for their matching data.^pField searching: The user is prompted with all the fields given in the dataset to select 2 for searching. The software uses the selected fields are variables and conducts a similar process to the search function, instead using the 2 fields to find appropriate matches. Similar parameters will be applied, however the strain is likely lower, as the function is only searching for matching fields instead of datatypes. Storing the results in an array again will help minimise this strain.^pGraphing function: Using the dropdown option to open the graphing screen, the user is prompted to select a field from the dataset for each axis to then use for graphing. The software takes the selected fields, and using imported graph capabilities, assigns the data from the fields to arrays for each axis, then plots the data and outputs to the user. The title given to the graph is synthesised during this process by filling in the blank name ‘<x-axis> vs <y-axis>’ with the respective axis headings. To minimise the software strain, the imports will only be called when the user is in the graphing screen.", 2)
    if (-not $ok) {
        Write-Host "ERROR: Find/Replace for 'Functions' heading failed"
    }

    # --- Step 4: fix up paragraph/list formatting on the newly-inserted paragraphs.
    #     They were all created inheriting the "Functions" heading's own Heading3 + outline
    #     numbering; reset that per the target layout. ---

    # Intro sentence -> plain body paragraph (no heading style, no list numbering)
    $pIntro = $d.Paragraphs.Item($targetIndex + 1)
    $pIntro.Style = "Normal"
    $pIntro.Range.ListFormat.ListType = 0

    # Three bullet points -> ListParagraph style, rejoin the existing numId=7 bullet list
    for ($off = 2; $off -le 4; $off++) {
        $pBullet = $d.Paragraphs.Item($targetIndex + $off)
        $pBullet.Style = "List Paragraph"
        $pBullet.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, 1)
    }

    # Three narrative paragraphs -> plain body paragraphs
    for ($off = 5; $off -le 7; $off++) {
        $pBody = $d.Paragraphs.Item($targetIndex + $off)
        $pBody.Style = "Normal"
        $pBody.Range.ListFormat.ListType = 0
    }

    Write-Host "Inserted 7 paragraphs after 'Functions' heading (paragraph $targetIndex); document now has $($d.Paragraphs.Count) paragraphs."
}
